# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Each entry updates the "Price" (D) and/or "Volume(1h)" (E) cell for a coin row.
# ForceText keeps values that look like plain numbers (e.g. "1.00", "7.74") stored
# as text, matching the sheet's existing inline-string convention for the Price column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Col = "D"; Value = "66.629.66"; ForceText = $false },
    @{ Row = 2; Col = "E"; Value = "  +0.21%  "; ForceText = $false },
    @{ Row = 3; Col = "D"; Value = "3.523.99"; ForceText = $false },
    @{ Row = 3; Col = "E"; Value = "  -1.91%  "; ForceText = $false },
    @{ Row = 4; Col = "D"; Value = "1.00"; ForceText = $true },
    @{ Row = 4; Col = "E"; Value = "  -0.09%  "; ForceText = $false },
    @{ Row = 5; Col = "D"; Value = "607.26"; ForceText = $true },
    @{ Row = 5; Col = "E"; Value = "  -0.23%  "; ForceText = $false },
    @{ Row = 6; Col = "D"; Value = "143.39"; ForceText = $true },
    @{ Row = 6; Col = "E"; Value = "  -3.80%  "; ForceText = $false },
    @{ Row = 7; Col = "D"; Value = "3.521.78"; ForceText = $false },
    @{ Row = 7; Col = "E"; Value = "  -1.97%  "; ForceText = $false },
    @{ Row = 8; Col = "D"; Value = "1.00"; ForceText = $true },
    @{ Row = 8; Col = "E"; Value = "  -0.23%  "; ForceText = $false },
    @{ Row = 9; Col = "D"; Value = "0.511"; ForceText = $true },
    @{ Row = 9; Col = "E"; Value = "  +4.30%  "; ForceText = $false },
    @{ Row = 10; Col = "D"; Value = "7.74"; ForceText = $true },
    @{ Row = 10; Col = "E"; Value = "  -3.57%  "; ForceText = $false },
    @{ Row = 11; Col = "E"; Value = "  -4.31%  "; ForceText = $false },
    @{ Row = 12; Col = "E"; Value = "  -1.59%  "; ForceText = $false },
    @{ Row = 13; Col = "D"; Value = "4.119.46"; ForceText = $false },
    @{ Row = 13; Col = "E"; Value = "  -2.12%  "; ForceText = $false },
    @{ Row = 14; Col = "E"; Value = "  -6.38%  "; ForceText = $false },
    @{ Row = 15; Col = "E"; Value = "  -3.76%  "; ForceText = $false },
    @{ Row = 16; Col = "D"; Value = "3.523.15"; ForceText = $false },
    @{ Row = 16; Col = "E"; Value = "  -1.29%  "; ForceText = $false },
    @{ Row = 17; Col = "E"; Value = "  +0.49%  "; ForceText = $false },
    @{ Row = 18; Col = "D"; Value = "66.486.79"; ForceText = $false },
    @{ Row = 19; Col = "D"; Value = "10.76"; ForceText = $true },
    @{ Row = 19; Col = "E"; Value = "  -6.93%  "; ForceText = $false },
    @{ Row = 20; Col = "E"; Value = "  -3.24%  "; ForceText = $false },
    @{ Row = 21; Col = "D"; Value = "14.68"; ForceText = $true },
    @{ Row = 21; Col = "E"; Value = "  -2.85%  "; ForceText = $false },
    @{ Row = 22; Col = "D"; Value = "423.23"; ForceText = $true },
    @{ Row = 22; Col = "E"; Value = "  -0.96%  "; ForceText = $false },
    @{ Row = 23; Col = "E"; Value = "  -4.88%  "; ForceText = $false },
    @{ Row = 24; Col = "E"; Value = "  -2.40%  "; ForceText = $false },
    @{ Row = 25; Col = "D"; Value = "3.662.76"; ForceText = $false },
    @{ Row = 25; Col = "E"; Value = "  -2.16%  "; ForceText = $false },
    @{ Row = 26; Col = "E"; Value = "  +0.14%  "; ForceText = $false },
    @{ Row = 27; Col = "E"; Value = "  -5.72%  "; ForceText = $false },
    @{ Row = 28; Col = "D"; Value = "7.96"; ForceText = $true },
    @{ Row = 28; Col = "E"; Value = "  -4.26%  "; ForceText = $false },
    @{ Row = 29; Col = "E"; Value = "  -1.99%  "; ForceText = $false },
    @{ Row = 30; Col = "D"; Value = "8.93"; ForceText = $true },
    @{ Row = 30; Col = "E"; Value = "  -5.18%  "; ForceText = $false },
    @{ Row = 31; Col = "D"; Value = "0.999"; ForceText = $true },
    @{ Row = 31; Col = "E"; Value = "  -0.13%  "; ForceText = $false },
    @{ Row = 32; Col = "D"; Value = "3.528.43"; ForceText = $false },
    @{ Row = 32; Col = "E"; Value = "  -1.89%  "; ForceText = $false },
    @{ Row = 33; Col = "D"; Value = "0.154"; ForceText = $true },
    @{ Row = 33; Col = "E"; Value = "  -2.02%  "; ForceText = $false },
    @{ Row = 34; Col = "D"; Value = "24.21"; ForceText = $true },
    @{ Row = 34; Col = "E"; Value = "  -4.86%  "; ForceText = $false },
    @{ Row = 36; Col = "D"; Value = "1.33"; ForceText = $true },
    @{ Row = 36; Col = "E"; Value = "  -9.22%  "; ForceText = $false },
    @{ Row = 37; Col = "D"; Value = "7.57"; ForceText = $true },
    @{ Row = 37; Col = "E"; Value = "  -3.57%  "; ForceText = $false },
    @{ Row = 38; Col = "E"; Value = "  -3.92%  "; ForceText = $false },
    @{ Row = 39; Col = "D"; Value = "173.64"; ForceText = $true },
    @{ Row = 39; Col = "E"; Value = "  -2.05%  "; ForceText = $false },
    @{ Row = 40; Col = "D"; Value = "5.22"; ForceText = $true },
    @{ Row = 40; Col = "E"; Value = "  -7.47%  "; ForceText = $false },
    @{ Row = 41; Col = "D"; Value = "0.0815"; ForceText = $true },
    @{ Row = 41; Col = "E"; Value = "  -4.87%  "; ForceText = $false },
    @{ Row = 42; Col = "D"; Value = "5.00"; ForceText = $true },
    @{ Row = 42; Col = "E"; Value = "  -4.43%  "; ForceText = $false },
    @{ Row = 43; Col = "D"; Value = "0.854"; ForceText = $true },
    @{ Row = 43; Col = "E"; Value = "  -4.65%  "; ForceText = $false },
    @{ Row = 44; Col = "D"; Value = "45.45"; ForceText = $true },
    @{ Row = 44; Col = "E"; Value = "  -0.99%  "; ForceText = $false },
    @{ Row = 45; Col = "E"; Value = "  -6.46%  "; ForceText = $false },
    @{ Row = 46; Col = "E"; Value = "  -0.03%  "; ForceText = $false },
    @{ Row = 47; Col = "D"; Value = "2.37"; ForceText = $true },
    @{ Row = 47; Col = "E"; Value = "  -7.85%  "; ForceText = $false },
    @{ Row = 48; Col = "D"; Value = "7.08"; ForceText = $true },
    @{ Row = 48; Col = "E"; Value = "  -1.77%  "; ForceText = $false },
    @{ Row = 49; Col = "D"; Value = "1.13"; ForceText = $true },
    @{ Row = 49; Col = "E"; Value = "  -4.13%  "; ForceText = $false },
    @{ Row = 50; Col = "D"; Value = "22.83"; ForceText = $true },
    @{ Row = 50; Col = "E"; Value = "  -4.88%  "; ForceText = $false },
    @{ Row = 51; Col = "D"; Value = "0.906"; ForceText = $true },
    @{ Row = 51; Col = "E"; Value = "  -4.92%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range("$($u.Col)$($u.Row)")
    if ($u.ForceText) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}
